$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cyclic rotation of weekly data across rows 3, 5 and 7:
#   new row3 <- old row5
#   new row5 <- old row7
#   new row7 <- old row3

$ws.Range("D3").Value = 44172
$ws.Range("M3").Value = 90
$ws.Range("N3").Value = 8500
$ws.Range("O3").Value = 9000
$ws.Range("P3").Value = 8806
$ws.Range("S3").Value = 629

$ws.Range("D5").Value = 44232
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11583
$ws.Range("S5").Value = 827

$ws.Range("D7").Value = 44229
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11364
$ws.Range("S7").Value = 812
